# Applies targeted value changes to the "Crédito disponível - Centraliza" sheet
# as described by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 9
$ws.Range("K9").Value = 153380.3

# Row 10
$ws.Range("O10").Value = 125812.36

# Row 15
$ws.Range("N15").Value = 482.16
$ws.Range("O15").Value = 482.16

# Row 18
$ws.Range("M18").Value = 1064451.35

# Row 19
$ws.Range("N19").Value = 8497.29

# Row 23
$ws.Range("N23").Value = 106312
$ws.Range("O23").Value = 106312

# Row 24
$ws.Range("N24").Value = 10207.2
$ws.Range("O24").Value = 5667.92
